# Scheduled-runner market price refresh: updates computed price/profit
# columns (H..N) produced by an external market-board fetch across the
# ALC, ARM, BSM, CRP, GSM, LTW and WVR leve-profit sheets. Pure data
# values - no formulas/styles involved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 1258.4
$ws.Cells.Item(29, 10).Value = 2474.4
$ws.Cells.Item(29, 12).Value = 7423.200000000001
$ws.Cells.Item(29, 14).Value = -7985.200000000001

$ws.Cells.Item(38, 8).Value = 489.57144

$ws.Cells.Item(43, 8).Value = 2915.1667
$ws.Cells.Item(43, 9).Value = 2163.6667
$ws.Cells.Item(43, 10).Value = 3666.6667
$ws.Cells.Item(43, 11).Value = 2163.6667
$ws.Cells.Item(43, 12).Value = 3666.6667
$ws.Cells.Item(43, 13).Value = -2094.6667
$ws.Cells.Item(43, 14).Value = -3804.6667

$ws.Cells.Item(58, 8).Value = 1645.3125
$ws.Cells.Item(58, 10).Value = 2420
$ws.Cells.Item(58, 12).Value = 7260
$ws.Cells.Item(58, 14).Value = -7560

$ws.Cells.Item(87, 8).Value = 27670.588
$ws.Cells.Item(87, 10).Value = 27670.588
$ws.Cells.Item(87, 12).Value = 27670.588
$ws.Cells.Item(87, 14).Value = -30166.588

$ws.Cells.Item(90, 8).Value = 27670.588
$ws.Cells.Item(90, 10).Value = 27670.588
$ws.Cells.Item(90, 12).Value = 83011.764
$ws.Cells.Item(90, 14).Value = -95491.764

$ws.Cells.Item(92, 8).Value = 852.76
$ws.Cells.Item(92, 9).Value = 852.76
$ws.Cells.Item(92, 11).Value = 852.76
$ws.Cells.Item(92, 13).Value = 395.24

$ws.Cells.Item(94, 8).Value = 5933.75
$ws.Cells.Item(94, 9).Value = 5120.5
$ws.Cells.Item(94, 11).Value = 5120.5
$ws.Cells.Item(94, 13).Value = -4669.5

$ws.Cells.Item(98, 8).Value = 4117.303
$ws.Cells.Item(98, 9).Value = 2619.4827
$ws.Cells.Item(98, 10).Value = 14976.5
$ws.Cells.Item(98, 11).Value = 2619.4827
$ws.Cells.Item(98, 12).Value = 14976.5
$ws.Cells.Item(98, 13).Value = -1121.4827
$ws.Cells.Item(98, 14).Value = -17972.5

$ws.Cells.Item(100, 8).Value = 2624.4614
$ws.Cells.Item(100, 9).Value = 2261.5
$ws.Cells.Item(100, 10).Value = 3834.3333
$ws.Cells.Item(100, 11).Value = 2261.5
$ws.Cells.Item(100, 12).Value = 3834.3333
$ws.Cells.Item(100, 13).Value = -1720.5
$ws.Cells.Item(100, 14).Value = -4916.3333

$ws.Cells.Item(113, 8).Value = 2492.55
$ws.Cells.Item(113, 9).Value = 1928.8462
$ws.Cells.Item(113, 10).Value = 3539.4285
$ws.Cells.Item(113, 11).Value = 1928.8462
$ws.Cells.Item(113, 12).Value = 3539.4285
$ws.Cells.Item(113, 13).Value = 1325.1538
$ws.Cells.Item(113, 14).Value = -10047.4285

$ws.Cells.Item(116, 8).Value = 16669336
$ws.Cells.Item(116, 9).Value = 33335382
$ws.Cells.Item(116, 10).Value = 3289.5
$ws.Cells.Item(116, 11).Value = 33335382
$ws.Cells.Item(116, 12).Value = 3289.5
$ws.Cells.Item(116, 13).Value = -33331940
$ws.Cells.Item(116, 14).Value = -10173.5

$ws.Cells.Item(122, 8).Value = 4117.303
$ws.Cells.Item(122, 9).Value = 2619.4827
$ws.Cells.Item(122, 10).Value = 14976.5
$ws.Cells.Item(122, 11).Value = 7858.4481
$ws.Cells.Item(122, 12).Value = 44929.5
$ws.Cells.Item(122, 13).Value = -5408.4481
$ws.Cells.Item(122, 14).Value = -49829.5

$ws.Cells.Item(138, 8).Value = 2669.8167
$ws.Cells.Item(138, 9).Value = 1593
$ws.Cells.Item(138, 10).Value = 4285.0415
$ws.Cells.Item(138, 11).Value = 4779
$ws.Cells.Item(138, 12).Value = 12855.1245
$ws.Cells.Item(138, 13).Value = 361
$ws.Cells.Item(138, 14).Value = -23135.1245

$ws.Cells.Item(141, 8).Value = 5164.9062
$ws.Cells.Item(141, 9).Value = 2024.8889
$ws.Cells.Item(141, 10).Value = 22121
$ws.Cells.Item(141, 11).Value = 6074.6667
$ws.Cells.Item(141, 12).Value = 66363
$ws.Cells.Item(141, 13).Value = -894.6666999999998
$ws.Cells.Item(141, 14).Value = -76723

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10686.7
$ws.Cells.Item(32, 9).Value = 10302.821
$ws.Cells.Item(32, 11).Value = 10302.821
$ws.Cells.Item(32, 13).Value = -10015.821

$ws.Cells.Item(102, 8).Value = 502505
$ws.Cells.Item(102, 9).Value = 5010
$ws.Cells.Item(102, 10).Value = 1000000
$ws.Cells.Item(102, 11).Value = 5010
$ws.Cells.Item(102, 12).Value = 1000000
$ws.Cells.Item(102, 13).Value = -3388
$ws.Cells.Item(102, 14).Value = -1003244

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 19135.857
$ws.Cells.Item(20, 9).Value = 24061.373
$ws.Cells.Item(20, 10).Value = 2843.7693
$ws.Cells.Item(20, 11).Value = 24061.373
$ws.Cells.Item(20, 12).Value = 2843.7693
$ws.Cells.Item(20, 13).Value = -23814.373
$ws.Cells.Item(20, 14).Value = -3337.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1154.05
$ws.Cells.Item(16, 9).Value = 1016.9091
$ws.Cells.Item(16, 10).Value = 1321.6666
$ws.Cells.Item(16, 11).Value = 1016.9091
$ws.Cells.Item(16, 12).Value = 1321.6666
$ws.Cells.Item(16, 13).Value = -729.9091
$ws.Cells.Item(16, 14).Value = -1895.6666

$ws.Cells.Item(113, 8).Value = 1154.05
$ws.Cells.Item(113, 9).Value = 1016.9091
$ws.Cells.Item(113, 10).Value = 1321.6666
$ws.Cells.Item(113, 11).Value = 1016.9091
$ws.Cells.Item(113, 12).Value = 1321.6666
$ws.Cells.Item(113, 13).Value = 1153.0909
$ws.Cells.Item(113, 14).Value = -5661.6666

$ws.Cells.Item(141, 8).Value = 39537.75
$ws.Cells.Item(141, 10).Value = 39537.75
$ws.Cells.Item(141, 12).Value = 39537.75
$ws.Cells.Item(141, 14).Value = -49897.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6051.852
$ws.Cells.Item(70, 9).Value = 5467.1875
$ws.Cells.Item(70, 10).Value = 6902.273
$ws.Cells.Item(70, 11).Value = 5467.1875
$ws.Cells.Item(70, 12).Value = 6902.273
$ws.Cells.Item(70, 13).Value = -5197.1875
$ws.Cells.Item(70, 14).Value = -7442.273

$ws.Cells.Item(73, 8).Value = 6051.852
$ws.Cells.Item(73, 9).Value = 5467.1875
$ws.Cells.Item(73, 10).Value = 6902.273
$ws.Cells.Item(73, 11).Value = 5467.1875
$ws.Cells.Item(73, 12).Value = 6902.273
$ws.Cells.Item(73, 13).Value = -4531.1875
$ws.Cells.Item(73, 14).Value = -8774.273000000001

$ws.Cells.Item(107, 8).Value = 1027.55
$ws.Cells.Item(107, 9).Value = 1254.3636
$ws.Cells.Item(107, 10).Value = 750.3333
$ws.Cells.Item(107, 11).Value = 1254.3636
$ws.Cells.Item(107, 12).Value = 750.3333
$ws.Cells.Item(107, 13).Value = 665.6364000000001
$ws.Cells.Item(107, 14).Value = -4590.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1566.6666
$ws.Cells.Item(46, 10).Value = 1671.4286
$ws.Cells.Item(46, 12).Value = 1671.4286
$ws.Cells.Item(46, 14).Value = -2047.4286

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 63746
$ws.Cells.Item(81, 9).Value = 54508.74
$ws.Cells.Item(81, 10).Value = 151500
$ws.Cells.Item(81, 11).Value = 109017.48
$ws.Cells.Item(81, 12).Value = 303000
$ws.Cells.Item(81, 13).Value = -107956.48
$ws.Cells.Item(81, 14).Value = -305122

$ws.Cells.Item(84, 8).Value = 63746
$ws.Cells.Item(84, 9).Value = 54508.74
$ws.Cells.Item(84, 10).Value = 151500
$ws.Cells.Item(84, 11).Value = 545087.4
$ws.Cells.Item(84, 12).Value = 1515000
$ws.Cells.Item(84, 13).Value = -539783.4
$ws.Cells.Item(84, 14).Value = -1525608

$ws.Cells.Item(96, 8).Value = 2198.3333
$ws.Cells.Item(96, 9).Value = 1720
$ws.Cells.Item(96, 11).Value = 1720
$ws.Cells.Item(96, 13).Value = -347

$ws.Cells.Item(100, 8).Value = 7093.8066
$ws.Cells.Item(100, 9).Value = 10957.685
$ws.Cells.Item(100, 10).Value = 976
$ws.Cells.Item(100, 11).Value = 21915.37
$ws.Cells.Item(100, 12).Value = 1952
$ws.Cells.Item(100, 13).Value = -21374.37
$ws.Cells.Item(100, 14).Value = -3034

$ws.Cells.Item(130, 8).Value = 23666.666
$ws.Cells.Item(130, 10).Value = 23666.666
$ws.Cells.Item(130, 12).Value = 23666.666
$ws.Cells.Item(130, 14).Value = -33706.666

$ws.Cells.Item(132, 8).Value = 1554.5227
$ws.Cells.Item(132, 9).Value = 951.5454999999999
$ws.Cells.Item(132, 11).Value = 2854.6365
$ws.Cells.Item(132, 13).Value = -324.6364999999996

$ws.Cells.Item(136, 8).Value = 1203.2858
$ws.Cells.Item(136, 9).Value = 1269.258
$ws.Cells.Item(136, 10).Value = 1017.36365
$ws.Cells.Item(136, 11).Value = 3807.774
$ws.Cells.Item(136, 12).Value = 3052.09095
$ws.Cells.Item(136, 13).Value = -1257.774
$ws.Cells.Item(136, 14).Value = -8152.09095
